# Apply edit described by the diff:
# On the "numeric" worksheet, set cell E2 to the text value "*"
# (matching the existing style already applied to that cell),
# and leave the selection positioned on E3 (as it would be after
# typing a value into E2 and pressing Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")
$ws.Activate()

$ws.Range("E2").Value = "*"

$ws.Range("E3").Select()
